# Apply "repull data, push all data, mean calculation" edit:
# Update the dSF column (F) values for the data rows to reflect the
# freshly re-pulled upstream data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    4  = -3
    5  = 8
    6  = 2
    7  = -5
    8  = 1
    9  = 6
    10 = -2
    11 = -2
    12 = 2
    13 = -2
    14 = 5
    15 = -3
    16 = 0
    17 = -1
    19 = 2
    20 = 5
    21 = -1
    23 = 1
    24 = -3
    25 = 3
    26 = -5
    27 = -6
    28 = 1
    29 = -2
    30 = 1
    31 = 7
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
